# Applies: "1. Added new (valve) msgs"
#
# Summary of the change:
#  - On sheet "event_log_type" two new valve-message rows are inserted
#    right before the existing "MSG_VALVE_EXEC_PATTERN" row (old row 50):
#      new row 50: MSG_VALVE_SET_CONT_MOVEMENT / trig_set_valve_cont /
#                  "valve seq movement set" / "valve seq movement fail"
#      new row 51: MSG_VALVE_SET_MOVE_DURATION / trig_valve_set_dur /
#                  "valve duration updated" / "valve duration set fail"
#    (all the old rows 50-53 shift down to 52-55)
#  - The active/selected sheet moves from "data_log_type" to
#    "event_log_type", whose selection becomes D51; "data_log_type"'s
#    selection becomes the whole of row 15 (A15:XFD15).

$wb = $excel.ActiveWorkbook

# --- Sheet: data_log_type ---
# loses tabSelected, selection moves from E4 to the entire row 15
$wsData = $wb.Worksheets.Item("data_log_type")
$wsData.Activate() | Out-Null
$wsData.Range("A15:XFD15").Select() | Out-Null

# --- Sheet: event_log_type ---
# gains tabSelected (becomes the active tab); two new rows are inserted
$wsEvent = $wb.Worksheets.Item("event_log_type")
$wsEvent.Activate() | Out-Null

# Insert a row at 50 and fill it in column order (A,B,C,D,E,F) with the
# MSG_VALVE_SET_MOVE_DURATION message - this row ends up as row 51 once
# the next row is inserted above it.
$wsEvent.Rows.Item(50).Insert() | Out-Null
$wsEvent.Cells.Item(50, 1).Value = "MSG_VALVE_SET_MOVE_DURATION"
$wsEvent.Cells.Item(50, 2).Value = 1
$wsEvent.Cells.Item(50, 3).Value = 1
$wsEvent.Cells.Item(50, 4).Value = "trig_valve_set_dur"
$wsEvent.Cells.Item(50, 5).Value = "valve duration updated"
$wsEvent.Cells.Item(50, 6).Value = "valve duration set fail"

# Insert another row above it (pushes the row above back down to 51) for
# MSG_VALVE_SET_CONT_MOVEMENT, filling D,E,F / B,C first and A last.
$wsEvent.Rows.Item(50).Insert() | Out-Null
$wsEvent.Cells.Item(50, 4).Value = "trig_set_valve_cont"
$wsEvent.Cells.Item(50, 5).Value = "valve seq movement set"
$wsEvent.Cells.Item(50, 6).Value = "valve seq movement fail"
$wsEvent.Cells.Item(50, 2).Value = 1
$wsEvent.Cells.Item(50, 3).Value = 1
$wsEvent.Cells.Item(50, 1).Value = "MSG_VALVE_SET_CONT_MOVEMENT"

# Selection on this sheet moves to D51
$wsEvent.Range("D51").Select() | Out-Null
